$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# Grab a couple of style "templates" that already exist in the sheet before
# anything else changes:
#   - the bold header style used by row 1
#   - the wrap-text style used by the old column A data cell (A2)
#   - the date-number-format style used by the old column B/C data cells
# We copy these onto the new target cells further down so every style index
# Excel ends up writing matches what a human editing the template by hand
# would produce (no stray/duplicate style entries).
# ---------------------------------------------------------------------------
$ws.Range("A1").Copy()
$ws.Range("A1:D1").PasteSpecial(-4122)  # xlPasteFormats (keep header style on all 4, incl. new A1)
$excel.CutCopyMode = $false

$ws.Range("B2").Copy()      # old "date" style (numFmt 14 + wrap + vcenter)
$ws.Range("C2:D5").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = $false

$ws.Range("A2").Copy()      # old "wrap text" style
$ws.Range("B2:B5").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = $false

$ws.Range("A2:A5").ClearFormats()   # new column A carries no special style
$excel.CutCopyMode = $false

# ---------------------------------------------------------------------------
# Column widths: new column A gets the narrower width, column B takes over
# the old column A width (37). Columns C/D are left completely untouched so
# they remain byte-for-byte identical to the original 30.6640625 width.
# ---------------------------------------------------------------------------
$ws.Columns("A").ColumnWidth = 26.75
$ws.Columns("B").ColumnWidth = 36.16

# ---------------------------------------------------------------------------
# Header labels: a new "Nombre_sede" column is introduced in front of the
# existing Descripcion / Fecha_inicio / Fecha_fin headers, which shift right
# by one. The old column D header ("Estado") is dropped entirely.
# ---------------------------------------------------------------------------
$ws.Range("D1").Value = $ws.Range("C1").Value2
$ws.Range("C1").Value = $ws.Range("B1").Value2
$ws.Range("B1").Value = $ws.Range("A1").Value2
$ws.Range("A1").Value = "Nombre_sede"

# ---------------------------------------------------------------------------
# Data values (4 rows, one per university campus).
# ---------------------------------------------------------------------------
$ws.Range("A2").Value = "Sede Colombia"
$ws.Range("B2").Value = "Proceso Ejemplo sede U colombia"
$ws.Range("C2").Value = 45108
$ws.Range("D2").Value = 45117

$ws.Range("A3").Value = "Sede Argentina"
$ws.Range("B3").Value = "Proceso Ejemplo sede U argentina"
$ws.Range("C3").Value = 45109
$ws.Range("D3").Value = 45122

$ws.Range("A4").Value = "Sede Venezuela"
$ws.Range("B4").Value = "Proceso Ejemplo sede U venezuela"
$ws.Range("C4").Value = 45110
$ws.Range("D4").Value = 45127

$ws.Range("A5").Value = "Sede Bélgica"
$ws.Range("B5").Value = "Proceso Ejemplo sede U bélgica"
$ws.Range("C5").Value = 45111
$ws.Range("D5").Value = 45132

$ws.Range("E8").Select()
